$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking values
# (e.g. "0.9986") are preserved as literal text/strings rather than
# being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.398.33"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.847.13"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "240.61"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "0.6343"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.07568"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "0.2968"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "2.887.27"
$ws.Range("E10").Value = "  +56.19%  "
$ws.Range("D11").Value = "24.56"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "0.07717"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "4.985"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "0.6849"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "82.76"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.000009932"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "6.189"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "29.423.82"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "231.73"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "7.587"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "154.71"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").Value = "0.1393"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").Value = "8.428"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "17.68"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "1.469"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "0.05805"
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "4.121"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("B33").Value = "RocketPoolETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D33").Value = "3.032.22"
$ws.Range("E33").Value = "  +51.49%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.870"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "1.160"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "0.7211"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "2.595"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "1.250.77"
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "2.790"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").Value = "0.01807"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "0.9038"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "6.069"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "0.9990"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "67.14"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("D46").Value = "7.318"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "9.144"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Value = "0.4012"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "1.696"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "0.05742"
$ws.Range("E51").Value = "  +0.05%  "

# Restore the original (default/general) cell style for column D so no
# extraneous per-cell style reference is left behind.
$ws.Range("D2:D51").Style = "Normal"

